$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed Price (D) / Volume(1h) (E) values scraped for this run.
# A handful of Price values are plain decimals (e.g. "307.57"); Excel would
# otherwise auto-convert those into numbers when typed into a General-
# formatted cell, so they're written with a leading apostrophe (quote-prefix)
# to keep them stored as text, matching the rest of the Price column.

$ws.Range('D2').Value = '27.217.76'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.904.94'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '''307.57'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '''0.5260'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').Value = '''0.3808'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = '''0.07308'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '''21.63'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').Value = '''0.9052'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '''0.08087'
$ws.Range('E12').Value = '  -3.97%  '
$ws.Range('D13').Value = '''95.78'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').Value = '''5.367'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '1.791.40'
$ws.Range('E15').Value = '  -5.96%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '''0.000008683'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '27.256.05'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '''5.129'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').Value = '''6.482'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').Value = '''2.357'
$ws.Range('E24').Value = '  +2.94%  '
$ws.Range('D25').Value = '''149.56'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').Value = '''18.27'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').Value = '''1.742'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').Value = '''117.23'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('D30').Value = '''4.890'
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').Value = '''0.09242'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = '''0.8043'
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').Value = '''0.05070'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = '''1.229'
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('D35').Value = '''2.977'
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('D36').Value = '''3.394'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').Value = '''2.697'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('D38').Value = '''0.5722'
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = '''0.01992'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').Value = '''1.086'
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('D41').Value = '''9.004'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '''6.604'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').Value = '''116.62'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').Value = '''0.1517'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = '''0.4916'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('D47').Value = '''10.15'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = '''1.642'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').Value = '''38.57'
$ws.Range('E49').Value = '  +2.97%  '
$ws.Range('D50').Value = '''64.30'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '''0.05964'
$ws.Range('E51').Value = '  +0.40%  '
